$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - weekday headers
$ws.Range("B1").Value = "segunda"
$ws.Range("C1").Value = "terça"
$ws.Range("D1").Value = "quarta"
$ws.Range("E1").Value = "quinta"
$ws.Range("F1").Value = "sexta"

# Row 2 - 7:00
$ws.Range("A2").Value = "7:00"
$ws.Range("C2").Value = "Desenho Técnico"

# Row 3 - 7:50
$ws.Range("A3").Value = "7:50"
$ws.Range("C3").Value = "EAP"

# Row 4 - 8:40
$ws.Range("A4").Value = "8:40"
$ws.Range("E4").Value = "-"

# Row 5 - 9:30
$ws.Range("A5").Value = "9:30"

# Row 6 - 10:40
$ws.Range("A6").Value = "10:40"

# Row 7 - 11:30
$ws.Range("A7").Value = "11:30"
$ws.Range("C7").Value = "Circuitos Elétricos 2"

# Row 8 - 13:00
$ws.Range("A8").Value = "13:00"
$ws.Range("F8").Value = "-"

# Row 9 - 13:50
$ws.Range("A9").Value = "13:50"

# Row 10 - 14:40
$ws.Range("A10").Value = "14:40"
$ws.Range("D10").Value = "EAP"

# Row 11 - 15:30
$ws.Range("A11").Value = "15:30"
$ws.Range("D11").Value = "-"

# Row 12 - 16:40
$ws.Range("A12").Value = "16:40"

# Row 13 - 17:30
$ws.Range("A13").Value = "17:30"
$ws.Range("F13").Value = "-"
